$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend dimension/used-range by copying the row-68 template formatting
# (this carries the bold/bordered style on column A and the date-time
# number format on column E onto the three newly appended rows).
$ws.Range("A68:V68").Copy() | Out-Null
$ws.Range("A69:V71").PasteSpecial(-4122) | Out-Null

# Row 69 (match index 68)
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = "poland"
$ws.Range("C69").Value = "ekstraklasa"
$ws.Range("D69").Value = "2023-2024"
$ws.Range("E69").Value = 45191.75
$ws.Range("F69").Value = "Korona Kielce"
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = "Widzew Lodz"
$ws.Range("I69").Value = 1
$ws.Range("J69").Value = 2.26
$ws.Range("K69").Value = "17/09/2023 14:13"
$ws.Range("L69").Value = 2.43
$ws.Range("M69").Value = "22/09/2023 17:51"
$ws.Range("N69").Value = 3.37
$ws.Range("O69").Value = "17/09/2023 14:13"
$ws.Range("P69").Value = 3.38
$ws.Range("Q69").Value = "22/09/2023 17:51"
$ws.Range("R69").Value = 3.1
$ws.Range("S69").Value = "17/09/2023 14:13"
$ws.Range("T69").Value = 3.07
$ws.Range("U69").Value = "22/09/2023 17:51"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/korona-kielce-widzew-lodz/Wxt0ATg3/"

# Row 70 (match index 69)
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = "poland"
$ws.Range("C70").Value = "ekstraklasa"
$ws.Range("D70").Value = "2023-2024"
$ws.Range("E70").Value = 45191.85416666666
$ws.Range("F70").Value = "LKS Lodz"
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = "Jagiellonia"
$ws.Range("I70").Value = 2
$ws.Range("J70").Value = 2.46
$ws.Range("K70").Value = "19/09/2023 13:42"
$ws.Range("L70").Value = 2.79
$ws.Range("M70").Value = "22/09/2023 20:25"
$ws.Range("N70").Value = 3.32
$ws.Range("O70").Value = "19/09/2023 13:42"
$ws.Range("P70").Value = 3.3
$ws.Range("Q70").Value = "22/09/2023 20:25"
$ws.Range("R70").Value = 2.84
$ws.Range("S70").Value = "19/09/2023 13:42"
$ws.Range("T70").Value = 2.7
$ws.Range("U70").Value = "22/09/2023 20:25"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/lks-lodz-jagiellonia/jsPr251q/"

# Row 71 (match index 70)
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "poland"
$ws.Range("C71").Value = "ekstraklasa"
$ws.Range("D71").Value = "2023-2024"
$ws.Range("E71").Value = 45192.52083333334
$ws.Range("F71").Value = "Radomiak Radom"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Puszcza"
$ws.Range("I71").Value = 1
$ws.Range("J71").Value = 1.59
$ws.Range("K71").Value = "19/09/2023 13:42"
$ws.Range("L71").Value = 1.6
$ws.Range("M71").Value = "23/09/2023 12:25"
$ws.Range("N71").Value = 4.08
$ws.Range("O71").Value = "19/09/2023 13:42"
$ws.Range("P71").Value = 4.19
$ws.Range("Q71").Value = "23/09/2023 12:25"
$ws.Range("R71").Value = 5.79
$ws.Range("S71").Value = "19/09/2023 13:42"
$ws.Range("T71").Value = 5.72
$ws.Range("U71").Value = "23/09/2023 12:25"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/radomiak-radom-puszcza/KbYC7RwM/"
